$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.364.25'
$ws.Range('E2').Value = '  +8.90%  '
$ws.Range('D3').Value = '1.599.76'
$ws.Range('E3').Value = '  +8.50%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D5').Value = '0.9939'
$ws.Range('E5').Value = '  +4.13%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '303.79'
$ws.Range('E6').Value = '  +9.37%  '
$ws.Range('D7').Value = '0.3654'
$ws.Range('E7').Value = '  +1.02%  '
$ws.Range('D8').Value = '0.3387'
$ws.Range('E8').Value = '  +10.53%  '
$ws.Range('D9').Value = '41.84'
$ws.Range('E9').Value = '  +6.04%  '
$ws.Range('D10').Value = '1.127'
$ws.Range('E10').Value = '  +6.02%  '
$ws.Range('D11').Value = '0.07025'
$ws.Range('E11').Value = '  +5.58%  '
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '19.63'
$ws.Range('E13').Value = '  +8.24%  '
$ws.Range('D14').Value = '5.885'
$ws.Range('E14').Value = '  +6.51%  '
$ws.Range('D15').Value = '6.601'
$ws.Range('E15').Value = '  +6.61%  '
$ws.Range('B16').Value = 'Dai'
$ws.Range('C16').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D16').Value = '0.9945'
$ws.Range('E16').Value = '  +4.21%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.00001074'
$ws.Range('E17').Value = '  +4.42%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '1.598.17'
$ws.Range('E18').Value = '  +8.34%  '
$ws.Range('D19').Value = '0.06600'
$ws.Range('E19').Value = '  +11.35%  '
$ws.Range('D20').Value = '77.25'
$ws.Range('E20').Value = '  +11.74%  '
$ws.Range('D21').Value = '5.991'
$ws.Range('E21').Value = '  +8.80%  '
$ws.Range('D22').Value = '15.89'
$ws.Range('E22').Value = '  +9.59%  '
$ws.Range('D23').Value = '11.74'
$ws.Range('E23').Value = '  +5.10%  '
$ws.Range('D24').Value = '22.385.22'
$ws.Range('E24').Value = '  +8.83%  '
$ws.Range('D25').Value = '2.383'
$ws.Range('E25').Value = '  +5.39%  '
$ws.Range('D26').Value = '2.541'
$ws.Range('E26').Value = '  +18.76%  '
$ws.Range('D27').Value = '149.12'
$ws.Range('E27').Value = '  +4.18%  '
$ws.Range('D28').Value = '19.41'
$ws.Range('E28').Value = '  +12.89%  '
$ws.Range('D29').Value = '1.778.00'
$ws.Range('E29').Value = '  +8.65%  '
$ws.Range('D30').Value = '122.00'
$ws.Range('E30').Value = '  +7.32%  '
$ws.Range('D31').Value = '4.046'
$ws.Range('E31').Value = '  +3.01%  '
$ws.Range('D32').Value = '6.086'
$ws.Range('E32').Value = '  +22.05%  '
$ws.Range('D33').Value = '0.9359'
$ws.Range('E33').Value = '  +15.95%  '
$ws.Range('D34').Value = '1.687'
$ws.Range('E34').Value = '  +11.48%  '
$ws.Range('D35').Value = '0.08190'
$ws.Range('E35').Value = '  +2.54%  '
$ws.Range('E36').Value = '  +14.37%  '
$ws.Range('E37').Value = '  +9.59%  '
$ws.Range('D38').Value = '1.249'
$ws.Range('E38').Value = '  +2.35%  '
$ws.Range('D39').Value = '8.470'
$ws.Range('E39').Value = '  +14.17%  '
$ws.Range('D40').Value = '0.06036'
$ws.Range('E40').Value = '  +3.20%  '
$ws.Range('E41').Value = '  +7.17%  '
$ws.Range('D42').Value = '0.2005'
$ws.Range('E42').Value = '  +6.77%  '
$ws.Range('D43').Value = '0.9943'
$ws.Range('E43').Value = '  +4.08%  '
$ws.Range('D44').Value = '0.5860'
$ws.Range('E44').Value = '  +10.76%  '
$ws.Range('D45').Value = '3.818'
$ws.Range('E45').Value = '  +8.31%  '
$ws.Range('D46').Value = '13.08'
$ws.Range('E46').Value = '  +6.53%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5638'
$ws.Range('E47').Value = '  +8.42%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '126.43'
$ws.Range('E48').Value = '  +6.76%  '
$ws.Range('D49').Value = '1.954'
$ws.Range('E49').Value = '  +7.60%  '
$ws.Range('D50').Value = '0.06776'
$ws.Range('E50').Value = '  +4.63%  '
$ws.Range('D51').Value = '73.09'
$ws.Range('E51').Value = '  +8.35%  '
